# PROS-13075 - CCRU - POS 2020 KPIs
#
# Insert a new "Sub brand" column right after the existing "Brand" column
# (column T) on the "HoReCa Bar Tavern_Night Club" sheet. Excel's native
# Insert-Column behaviour shifts every column from U onward one slot to
# the right and has the new column inherit formatting from the column
# immediately to its left (old column T / "Brand"), which is exactly what
# we replicate here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HoReCa Bar Tavern_Night Club")
$ws.Activate()

# Insert a new blank column at U (21st column); everything from U onward
# (including the "Logical Operator" header that used to live in U) shifts
# one column to the right automatically.
$ws.Columns("U:U").Insert()

# Name the newly inserted header cell.
$ws.Range("U1").Value = "Sub brand"

# Re-establish the AutoFilter over the now-wider used range (A1:AP31).
# Toggling it off first avoids Excel's "re-apply same range" no-op/removal
# quirk when the filter already covers (part of) the target range.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:AP31").AutoFilter()

# Keep the workbook-level "_FilterDatabase" defined names (hidden xlnm one
# plus the visible helper copy) in sync with the widened filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='HoReCa Bar Tavern_Night Club'!`$A`$1:`$AP`$31"
    }
}

# Reflect the author's final cursor position after typing the new header
# and checking the bottom of the new column.
[void]$ws.Range("U31").Select()
